# Generate Report for Handback
# The dc34a600-...md file has now been handed back (in sync with en-US),
# so every sheet's row for that file needs its status/time refreshed and
# its stale "version mismatch" error message cleared.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-22 14:55:18"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 12.85

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("K2").Value = "2016-08-22 14:54:22"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-22 14:55:27"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 12.85
